$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Checklist")

# Insert a new row at position 70 (everything from row 70 downward shifts down by one).
$ws.Range("A70:G70").Insert(-4121)

# Seed the new row's formatting (borders/fonts/alignment/row height) by copying the
# immediately preceding row (69, "WSTG-SESS-10"), which has matching column layout.
$ws.Range("A69:G69").Copy($ws.Range("A70:G70"))
$ws.Rows.Item(70).RowHeight = 49.5

# Now overwrite the copied content with the new "WSTG-SESS-11" entry.
$ws.Range("A70").Value = $null
$ws.Range("B70").Value = "WSTG-SESS-11"
$ws.Range("C70").Formula = '=HYPERLINK("https://owasp.org/www-project-web-security-testing-guide/latest/4-Web_Application_Security_Testing/06-Session_Management_Testing/11-Testing_for_Concurrent_Sessions", "Testing for Concurrent Sessions")'
$ws.Range("D70").Value = "- Evaluate the application's session management by assessing the handling of multiple active sessions for a single user account."
$ws.Range("E70").Value = "Not Started"
$ws.Range("F70").Value = $null

Write-Output "done"
